$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.UsedRange.ClearContents()

$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "peso"
$ws.Cells.Item(1,4).Value = "tipo"
$ws.Cells.Item(1,5).Value = "price"
$ws.Cells.Item(1,6).Value = "precioSuelto"
$ws.Cells.Item(1,7).Value = "stock"
$ws.Cells.Item(1,8).Value = "fechaActualizacion"
$ws.Cells.Item(1,9).Value = "tipoEdad"
$ws.Cells.Item(1,10).Value = "tipoMordida"
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "DOGUI"
$ws.Cells.Item(2,3).Value = 23
$ws.Cells.Item(2,4).Value = "PERRO"
$ws.Cells.Item(2,5).Value = 23
$ws.Cells.Item(2,6).Value = 20
$ws.Cells.Item(2,7).Value = 23
$ws.Cells.Item(2,8).Value = "'2023-09-08"
$ws.Cells.Item(2,8).ClearFormats()
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "NUTRIBON"
$ws.Cells.Item(3,3).Value = "'20"
$ws.Cells.Item(3,3).ClearFormats()
$ws.Cells.Item(3,4).Value = "PERRO"
$ws.Cells.Item(3,5).Value = "'31"
$ws.Cells.Item(3,5).ClearFormats()
$ws.Cells.Item(3,7).Value = "'31"
$ws.Cells.Item(3,7).ClearFormats()
$ws.Cells.Item(3,8).Value = "'2023-09-07"
$ws.Cells.Item(3,8).ClearFormats()
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "A"
$ws.Cells.Item(4,3).Value = "'123"
$ws.Cells.Item(4,3).ClearFormats()
$ws.Cells.Item(4,4).Value = "PERRO"
$ws.Cells.Item(4,5).Value = "'123"
$ws.Cells.Item(4,5).ClearFormats()
$ws.Cells.Item(4,7).Value = "'123"
$ws.Cells.Item(4,7).ClearFormats()
$ws.Cells.Item(4,8).Value = "'2023-09-07"
$ws.Cells.Item(4,8).ClearFormats()
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "B"
$ws.Cells.Item(5,3).Value = "'123"
$ws.Cells.Item(5,3).ClearFormats()
$ws.Cells.Item(5,4).Value = "PERRO"
$ws.Cells.Item(5,5).Value = "'123"
$ws.Cells.Item(5,5).ClearFormats()
$ws.Cells.Item(5,7).Value = "'123"
$ws.Cells.Item(5,7).ClearFormats()
$ws.Cells.Item(5,8).Value = "'2023-09-07"
$ws.Cells.Item(5,8).ClearFormats()
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "NUTRIBON"
$ws.Cells.Item(6,3).Value = "'8"
$ws.Cells.Item(6,3).ClearFormats()
$ws.Cells.Item(6,4).Value = "GATO"
$ws.Cells.Item(6,5).Value = "'10000"
$ws.Cells.Item(6,5).ClearFormats()
$ws.Cells.Item(6,7).Value = "'1"
$ws.Cells.Item(6,7).ClearFormats()
$ws.Cells.Item(6,8).Value = "'2023-09-07"
$ws.Cells.Item(6,8).ClearFormats()
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "NUTRIBON"
$ws.Cells.Item(7,3).Value = "'12"
$ws.Cells.Item(7,3).ClearFormats()
$ws.Cells.Item(7,4).Value = "PERRO"
$ws.Cells.Item(7,5).Value = "'123123"
$ws.Cells.Item(7,5).ClearFormats()
$ws.Cells.Item(7,7).Value = "'22"
$ws.Cells.Item(7,7).ClearFormats()
$ws.Cells.Item(7,8).Value = "'2023-09-07"
$ws.Cells.Item(7,8).ClearFormats()
$ws.Cells.Item(7,9).Value = "ADULTO"
$ws.Cells.Item(7,10).Value = "GRANDE"
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "JUAN"
$ws.Cells.Item(8,3).Value = "'2"
$ws.Cells.Item(8,3).ClearFormats()
$ws.Cells.Item(8,4).Value = "PERRO"
$ws.Cells.Item(8,5).Value = "'1777"
$ws.Cells.Item(8,5).ClearFormats()
$ws.Cells.Item(8,6).Value = "'23"
$ws.Cells.Item(8,6).ClearFormats()
$ws.Cells.Item(8,7).Value = "'1"
$ws.Cells.Item(8,7).ClearFormats()
$ws.Cells.Item(8,8).Value = "'2023-09-08"
$ws.Cells.Item(8,8).ClearFormats()
$ws.Cells.Item(8,9).Value = "CACHORRO"
$ws.Cells.Item(8,10).Value = "GRANDE"
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "ASD"
$ws.Cells.Item(9,3).Value = "'10"
$ws.Cells.Item(9,3).ClearFormats()
$ws.Cells.Item(9,4).Value = "PERRO"
$ws.Cells.Item(9,5).Value = "'123"
$ws.Cells.Item(9,5).ClearFormats()
$ws.Cells.Item(9,6).Value = "'123"
$ws.Cells.Item(9,6).ClearFormats()
$ws.Cells.Item(9,7).Value = "'123"
$ws.Cells.Item(9,7).ClearFormats()
$ws.Cells.Item(9,8).Value = "'2023-09-08"
$ws.Cells.Item(9,8).ClearFormats()
$ws.Cells.Item(9,9).Value = "-"
$ws.Cells.Item(9,10).Value = "-"
